$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: Id value change
$ws.Range("A2").Value = 89602968

# J2, K2, N2: clear these (previously empty inline-string) cells entirely
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("N2").ClearContents()

# S2: Noggrannhet value change
$ws.Range("S2").Value = 10

# AC2: remove the public comment cell entirely
$ws.Range("AC2").ClearContents()

# AF2: remove the (empty) determination-method cell entirely
$ws.Range("AF2").ClearContents()

# AX2: Observatörer text change
$ws.Range("AX2").Value = "Via Erland Lindblad"

# AY2: Projektnamn text set
$ws.Range("AY2").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"
